$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right after the header row (row 1) to make room for the
# new "Agosto 2025" (08/2025) daily entries. This pushes all the existing
# July/June/May data down by 3 rows, exactly matching the diff.
$ws.Rows("2:4").Insert()

# Fill the 3 new rows with the new August 2025 data (Dia, total_venda, Mes, Ano, Periodo)
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 21751.6
$ws.Cells.Item(2, 3).Value = 8
$ws.Cells.Item(2, 4).Value = 2025
$ws.Cells.Item(2, 5).Value = "08/2025"

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 20735.7
$ws.Cells.Item(3, 3).Value = 8
$ws.Cells.Item(3, 4).Value = 2025
$ws.Cells.Item(3, 5).Value = "08/2025"

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 5020
$ws.Cells.Item(4, 3).Value = 8
$ws.Cells.Item(4, 4).Value = 2025
$ws.Cells.Item(4, 5).Value = "08/2025"
